$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Values to set for column 4 (1077-width, "Diem lan 1") and column 5 (1198-width, "Diem lan 2")
# for data rows 2..11 (row 1 is the header).
$col4 = @{2="1"; 3="1"; 4="0"; 5="1"; 6="0"; 7="0"; 8="1"; 9="0"; 10="1"; 11="1"}
$col5 = @{2="1"; 3="1"; 4="1"; 5="1"; 6="1"; 7="1"; 8="1"; 9="0"; 10="1"; 11="1"}

foreach ($r in 2..11) {
    $row = $t.Rows.Item($r)

    $c4 = $row.Cells.Item(4).Range
    $c4.End = $c4.End - 1
    $c4.Text = $col4[$r]

    $c5 = $row.Cells.Item(5).Range
    $c5.End = $c5.End - 1
    $c5.Text = $col5[$r]
}

# Fill in the final (previously empty) paragraph after the table with
# nine tab stops followed by the "Tong" summary line. Locate it via the
# table's end position rather than Paragraphs.Last, since paragraph
# indexing inside/after the table can drift once cell text has changed.
# The range spans the whole paragraph (including its end mark) so that
# InsertXML replaces it cleanly.
$rng = $d.Range($t.Range.End, $d.Content.End)

$xmlFrag = "<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Tổng:      6</w:t></w:r><w:r><w:tab/><w:t xml:space=`"preserve`">          9</w:t></w:r></w:p>"
[void]$rng.InsertXML($xmlFrag)
